$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet "Sheet" -> "Sheet1"
$ws.Name = "Sheet1"

# 2. Insert a new column before column A to make room for the new "ID" column
#    (shifts horario/nome/email/empresa/telefone from A:E to B:F)
$ws.Columns.Item(1).Insert()

# 3. Set the new header cell
$a1 = $ws.Range("A1")
$a1.Value = "ID"

# 4. Style the new header cell: bold font, thin border all around,
#    centered horizontally, top-aligned vertically.
$a1.Font.Bold = $true
$a1.Borders.LineStyle = 1
$a1.HorizontalAlignment = -4108
$a1.VerticalAlignment = -4160

# 5. Apply the same formatting to the rest of the header row (B1:F1)
$a1.Copy()
$ws.Range("B1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 6. Add the new reservation row (row 2)
$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 2).Value = "13:30"
$ws.Cells.Item(2, 3).Value = "joao cleiton"
$ws.Cells.Item(2, 4).Value = "igor.carneiro@mt.senac.br"
$ws.Cells.Item(2, 5).Value = "teste"
$ws.Cells.Item(2, 6).Value = 65999196160
